# bug fix in Eduati data files
#
# Sheet1 ("COLO320HSR_noCTRL_meas") had 43 extra rows (45:87) left over from
# a previous paste that only carried a stray running index in column A with
# no real measurement data - remove them so the sheet matches the other two
# (A1:N44). Also refresh the view state (selection / active sheet) to where
# the workbook was left after the cleanup.

$wb = $excel.ActiveWorkbook

# --- Sheet1: drop the bogus tail rows (45:87 only ever held a leftover index) ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Rows("45:87").Delete() | Out-Null

# --- Sheet3: it was the tab left active/selected before the fix; update its
#     lingering selection and let Sheet1 take over as the active tab below ---
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("B43").Select() | Out-Null

# --- Sheet1 becomes the active sheet/tab, scrolled down with F62 selected ---
$ws1.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("F62").Select() | Out-Null
